$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.851.32"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.600.63"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.06%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.77"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +2.88%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.84"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.20%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +1.02%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.69"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D13").Value = "3.058.01"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "60.858.17"
$ws.Range("E14").Value = "  +0.45%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.64"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "2.602.58"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  -0.20%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "355.35"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +2.40%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.57"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +0.61%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.20"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +1.19%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +0.07%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.98"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "2.718.45"
$ws.Range("E26").Value = "  +0.58%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -0.32%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.39"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +0.03%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.29"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +10.13%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.41"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  +2.56%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.17"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("E35").Value = "  +4.65%  "
$ws.Range("E36").Value = "  +0.79%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.916"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +7.52%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.909"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +7.22%  "
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("E41").Value = "  +0.81%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "291.52"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("E43").Value = "  +1.52%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.622"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -0.18%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0559"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -0.14%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.59"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  +0.39%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.18"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +7.61%  "
